$d = $word.ActiveDocument

# Update the date line at the top of the document (unique text, safe to
# use a document-wide Find/Replace)
$d.Content.Find.Execute("2025-06-11 Wednesday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-06-12 Thursday", 2) | Out-Null

# Update the division problems in the table. Several cells share the same
# original text, so address each cell positionally (Table.Cell(row, col))
# and set its Range.Text directly rather than using Find/Replace, which
# would match the first occurrence anywhere in the document.
$t = $d.Tables.Item(1)

$cell = $t.Cell(1,1)
# was "70÷2="
$cell.Range.Text = "42÷9="

$cell = $t.Cell(1,2)
# was "76÷4="
$cell.Range.Text = "68÷3="

$cell = $t.Cell(1,3)
# was "79÷6="
$cell.Range.Text = "92÷7="

$cell = $t.Cell(1,4)
# was "49÷5="
$cell.Range.Text = "16÷4="

$cell = $t.Cell(1,5)
# was "41÷8="
$cell.Range.Text = "20÷7="

$cell = $t.Cell(5,1)
# was "66÷6="
$cell.Range.Text = "77÷4="

$cell = $t.Cell(5,2)
# was "82÷8="
$cell.Range.Text = "47÷5="

$cell = $t.Cell(5,3)
# was "67÷2="
$cell.Range.Text = "44÷6="

$cell = $t.Cell(5,4)
# was "64÷6="
$cell.Range.Text = "84÷8="

$cell = $t.Cell(5,5)
# was "38÷3="
$cell.Range.Text = "99÷8="

$cell = $t.Cell(9,1)
# was "82÷8="
$cell.Range.Text = "55÷9="

$cell = $t.Cell(9,2)
# was "33÷3="
$cell.Range.Text = "29÷6="

$cell = $t.Cell(9,3)
# was "68÷6="
$cell.Range.Text = "93÷7="

$cell = $t.Cell(9,4)
# was "82÷9="
$cell.Range.Text = "68÷3="

$cell = $t.Cell(9,5)
# was "59÷9="
$cell.Range.Text = "92÷3="

$cell = $t.Cell(13,1)
# was "46÷2="
$cell.Range.Text = "56÷2="

$cell = $t.Cell(13,2)
# was "25÷7="
$cell.Range.Text = "95÷9="

$cell = $t.Cell(13,3)
# was "61÷5="
$cell.Range.Text = "74÷5="

$cell = $t.Cell(13,4)
# was "22÷8="
$cell.Range.Text = "16÷7="

$cell = $t.Cell(13,5)
# was "27÷3="
$cell.Range.Text = "77÷2="

$cell = $t.Cell(17,1)
# was "97÷7="
$cell.Range.Text = "49÷5="

$cell = $t.Cell(17,2)
# was "47÷9="
$cell.Range.Text = "66÷2="

$cell = $t.Cell(17,3)
# was "79÷6="
$cell.Range.Text = "76÷4="

$cell = $t.Cell(17,4)
# was "20÷4="
$cell.Range.Text = "24÷5="

$cell = $t.Cell(17,5)
# was "67÷8="
$cell.Range.Text = "18÷8="
